# [#134279773] Importacao de usuarios aceita informacao de turma
# Adds a new "Turma" (class/cohort) column to the import template, with
# sample values for the first two sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column J
$ws.Range("J1").Value = "Turma"

# Sample turma values for the first two sample people; the third sample
# row (Luana) is left without a turma, same as the source data.
$ws.Range("J2").Value = "Turma A"
$ws.Range("J3").Value = "Turma B"

# Move/restore the active selection to the newly added cell, matching the
# state the workbook was left in after the edit.
$ws.Range("J4").Select() | Out-Null

# Cosmetic: slightly tighter sheet tab ratio in the saved workbook view.
try {
    $excel.ActiveWindow.TabRatio = 985 / 1650
} catch {
}
